$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44606
$ws.Cells.Item(2, 13).Value = 45
$ws.Cells.Item(2, 18).Value = 'Provincia de Linares'

$ws.Cells.Item(3, 4).Value = 44592
$ws.Cells.Item(3, 13).Value = 30
$ws.Cells.Item(3, 14).Value = 8000
$ws.Cells.Item(3, 15).Value = 8000
$ws.Cells.Item(3, 16).Value = 8000
$ws.Cells.Item(3, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(3, 19).Value = 4000

$ws.Cells.Item(4, 4).Value = 44614
$ws.Cells.Item(4, 13).Value = 45
$ws.Cells.Item(4, 14).Value = 6000
$ws.Cells.Item(4, 15).Value = 6000
$ws.Cells.Item(4, 16).Value = 6000
$ws.Cells.Item(4, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(4, 19).Value = 3000

$ws.Cells.Item(5, 4).Value = 44585
$ws.Cells.Item(5, 13).Value = 160
$ws.Cells.Item(5, 14).Value = 6500
$ws.Cells.Item(5, 15).Value = 7000
$ws.Cells.Item(5, 16).Value = 6750
$ws.Cells.Item(5, 19).Value = 3375

$ws.Cells.Item(6, 4).Value = 44974
$ws.Cells.Item(6, 13).Value = 130
$ws.Cells.Item(6, 14).Value = 7000
$ws.Cells.Item(6, 15).Value = 7500
$ws.Cells.Item(6, 16).Value = 7269
$ws.Cells.Item(6, 19).Value = 3634

$ws.Cells.Item(7, 4).Value = 44582
$ws.Cells.Item(7, 13).Value = 150
$ws.Cells.Item(7, 14).Value = 6000
$ws.Cells.Item(7, 15).Value = 6500
$ws.Cells.Item(7, 16).Value = 6233
$ws.Cells.Item(7, 19).Value = 3116

$ws.Cells.Item(8, 4).Value = 44211
$ws.Cells.Item(8, 13).Value = 45
$ws.Cells.Item(8, 14).Value = 6000
$ws.Cells.Item(8, 15).Value = 6000
$ws.Cells.Item(8, 16).Value = 6000
$ws.Cells.Item(8, 19).Value = 3000

$ws.Cells.Item(9, 4).Value = 44214
$ws.Cells.Item(9, 13).Value = 48

$ws.Cells.Item(10, 4).Value = 44586
$ws.Cells.Item(10, 13).Value = 80
$ws.Cells.Item(10, 14).Value = 7000
$ws.Cells.Item(10, 15).Value = 7000
$ws.Cells.Item(10, 16).Value = 7000
$ws.Cells.Item(10, 19).Value = 3500

$ws.Cells.Item(11, 4).Value = 44588
$ws.Cells.Item(11, 13).Value = 160
$ws.Cells.Item(11, 14).Value = 6500
$ws.Cells.Item(11, 16).Value = 6750
$ws.Cells.Item(11, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(11, 19).Value = 3375

$ws.Cells.Item(12, 4).Value = 44628
$ws.Cells.Item(12, 13).Value = 40
$ws.Cells.Item(12, 14).Value = 6000
$ws.Cells.Item(12, 15).Value = 6000
$ws.Cells.Item(12, 16).Value = 6000
$ws.Cells.Item(12, 19).Value = 3000

$ws.Cells.Item(13, 4).Value = 44627

$ws.Cells.Item(14, 4).Value = 44589
$ws.Cells.Item(14, 13).Value = 60
$ws.Cells.Item(14, 15).Value = 6000
$ws.Cells.Item(14, 16).Value = 6000
$ws.Cells.Item(14, 19).Value = 3000

$ws.Cells.Item(15, 4).Value = 44587
$ws.Cells.Item(15, 13).Value = 165
$ws.Cells.Item(15, 16).Value = 6742
$ws.Cells.Item(15, 18).Value = 'Provincia de Linares'
$ws.Cells.Item(15, 19).Value = 3371

$ws.Cells.Item(16, 4).Value = 44960
$ws.Cells.Item(16, 13).Value = 40
$ws.Cells.Item(16, 14).Value = 7000
$ws.Cells.Item(16, 15).Value = 7000
$ws.Cells.Item(16, 16).Value = 7000
$ws.Cells.Item(16, 19).Value = 3500

$ws.Cells.Item(17, 4).Value = 44209
$ws.Cells.Item(17, 13).Value = 58
$ws.Cells.Item(17, 18).Value = 'Provincia de Curicó'

$ws.Cells.Item(18, 4).Value = 44959
$ws.Cells.Item(18, 13).Value = 40
$ws.Cells.Item(18, 14).Value = 7000
$ws.Cells.Item(18, 15).Value = 7000
$ws.Cells.Item(18, 16).Value = 7000
$ws.Cells.Item(18, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(18, 19).Value = 3500

$ws.Cells.Item(19, 4).Value = 45001
$ws.Cells.Item(19, 13).Value = 66
$ws.Cells.Item(19, 14).Value = 7500
$ws.Cells.Item(19, 16).Value = 7773
$ws.Cells.Item(19, 18).Value = 'Provincia de Curicó'
$ws.Cells.Item(19, 19).Value = 3886
